$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that directly follows the
#    title heading (it is being relocated, reworded slightly, to the end of
#    the document).
$metaPara = $d.Paragraphs.Item(2)
$null = $metaPara.Range.Delete()

# 2. Insert a new paragraph right before the final (closing image-prompt)
#    paragraph that carries a bold "Play Fever for Free - Slot Game Review"
#    run - this mirrors the title that used to sit in the removed meta
#    description paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$null = $lastPara.Range.InsertParagraphBefore()

$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fever for Free - Slot Game Review</w:t></w:r></w:p>'
$newTitlePara = $d.Paragraphs.Item($count)
$null = $newTitlePara.Range.InsertXML($titleXml)

# 3. Replace the closing paragraph's italic image-prompt text with the meta
#    description copy that used to live near the top of the document.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$oldText = $finalPara.Range.Text
$null = $finalPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Fever by Cristaltec, a slot game that combines classic mechanics with modern features. Play for free and trigger free spins!", 2)
